$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 245
$ws.Range("I96").Value = 212.83333
$ws.Range("J96").Value = 341.5
$ws.Range("K96").Value = 638.49999
$ws.Range("L96").Value = 1024.5
$ws.Range("M96").Value = 734.50001
$ws.Range("N96").Value = -3770.5
$ws.Range("H111").Value = 13808.333
$ws.Range("I111").Value = 4548.1665
$ws.Range("J111").Value = 23068.5
$ws.Range("K111").Value = 13644.4995
$ws.Range("L111").Value = 69205.5
$ws.Range("M111").Value = -10577.4995
$ws.Range("N111").Value = -75339.5
$ws.Range("H137").Value = 1985.7693
$ws.Range("I137").Value = 1105.625
$ws.Range("K137").Value = 3316.875
$ws.Range("M137").Value = -766.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2945.1
$ws.Range("I61").Value = 2507.2856
$ws.Range("K61").Value = 2507.2856
$ws.Range("M61").Value = -2295.2856
$ws.Range("H74").Value = 1238.4
$ws.Range("I74").Value = 1450
$ws.Range("J74").Value = 979.7778
$ws.Range("K74").Value = 1450
$ws.Range("L74").Value = 979.7778
$ws.Range("M74").Value = -576
$ws.Range("N74").Value = -2727.7778
$ws.Range("H77").Value = 1238.4
$ws.Range("I77").Value = 1450
$ws.Range("J77").Value = 979.7778
$ws.Range("K77").Value = 7250
$ws.Range("L77").Value = 4898.889
$ws.Range("M77").Value = -2882
$ws.Range("N77").Value = -13634.889
$ws.Range("H122").Value = 7483
$ws.Range("I122").Value = 8596.799999999999
$ws.Range("J122").Value = 1914
$ws.Range("K122").Value = 25790.4
$ws.Range("L122").Value = 5742
$ws.Range("M122").Value = -23340.4
$ws.Range("N122").Value = -10642
$ws.Range("H136").Value = 2945.1
$ws.Range("I136").Value = 2507.2856
$ws.Range("K136").Value = 7521.8568
$ws.Range("M136").Value = -4971.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2780.6858
$ws.Range("I134").Value = 2348.8696
$ws.Range("J134").Value = 3608.3333
$ws.Range("K134").Value = 7046.6088
$ws.Range("L134").Value = 10824.9999
$ws.Range("M134").Value = -4511.6088
$ws.Range("N134").Value = -15894.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2116.561
$ws.Range("I31").Value = 1319.5588
$ws.Range("J31").Value = 5987.7144
$ws.Range("K31").Value = 1319.5588
$ws.Range("L31").Value = 5987.7144
$ws.Range("M31").Value = -1024.5588
$ws.Range("N31").Value = -6577.7144
$ws.Range("H34").Value = 2116.561
$ws.Range("I34").Value = 1319.5588
$ws.Range("J34").Value = 5987.7144
$ws.Range("K34").Value = 1319.5588
$ws.Range("L34").Value = 5987.7144
$ws.Range("M34").Value = -1117.5588
$ws.Range("N34").Value = -6391.7144
$ws.Range("H58").Value = 1456.0526
$ws.Range("I58").Value = 1103.4073
$ws.Range("K58").Value = 1103.4073
$ws.Range("M58").Value = -900.4073000000001
$ws.Range("H107").Value = 484.45
$ws.Range("I107").Value = 437
$ws.Range("J107").Value = 504.7857
$ws.Range("K107").Value = 437
$ws.Range("L107").Value = 504.7857
$ws.Range("M107").Value = 1483
$ws.Range("N107").Value = -4344.7857
$ws.Range("H132").Value = 3977.6924
$ws.Range("I132").Value = 2851.5
$ws.Range("J132").Value = 5779.6
$ws.Range("K132").Value = 8554.5
$ws.Range("L132").Value = 17338.8
$ws.Range("M132").Value = -6024.5
$ws.Range("N132").Value = -22398.8
$ws.Range("H134").Value = 2128.2593
$ws.Range("I134").Value = 1699.2632
$ws.Range("J134").Value = 3147.125
$ws.Range("K134").Value = 5097.7896
$ws.Range("L134").Value = 9441.375
$ws.Range("M134").Value = -2562.7896
$ws.Range("N134").Value = -14511.375
$ws.Range("H136").Value = 1456.0526
$ws.Range("I136").Value = 1103.4073
$ws.Range("K136").Value = 3310.2219
$ws.Range("M136").Value = -760.2219000000005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 32687.5
$ws.Range("J93").Value = 32687.5
$ws.Range("L93").Value = 32687.5
$ws.Range("N93").Value = -36431.5
$ws.Range("H122").Value = 4273.3335
$ws.Range("J122").Value = 4273.3335
$ws.Range("L122").Value = 12820.0005
$ws.Range("N122").Value = -17720.0005
$ws.Range("H123").Value = 8824.549999999999
$ws.Range("J123").Value = 8824.549999999999
$ws.Range("L123").Value = 8824.549999999999
$ws.Range("N123").Value = -13724.55
$ws.Range("H126").Value = 3428.3635
$ws.Range("I126").Value = 3366.1177
$ws.Range("J126").Value = 3640
$ws.Range("K126").Value = 10098.3531
$ws.Range("L126").Value = 10920
$ws.Range("M126").Value = -7628.3531
$ws.Range("N126").Value = -15860
$ws.Range("H131").Value = 25300
$ws.Range("J131").Value = 25300
$ws.Range("L131").Value = 25300
$ws.Range("N131").Value = -35380
$ws.Range("H133").Value = 59280
$ws.Range("J133").Value = 59280
$ws.Range("L133").Value = 59280
$ws.Range("N133").Value = -69400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 90002800
$ws.Range("I122").Value = 125002500
$ws.Range("K122").Value = 375007500
$ws.Range("M122").Value = -375005050
$ws.Range("H132").Value = 4793.3184
$ws.Range("I132").Value = 4663.5713
$ws.Range("J132").Value = 5020.375
$ws.Range("K132").Value = 13990.7139
$ws.Range("L132").Value = 15061.125
$ws.Range("M132").Value = -11460.7139
$ws.Range("N132").Value = -20121.125
$ws.Range("H136").Value = 2726.2
$ws.Range("I136").Value = 2010.3334
$ws.Range("J136").Value = 3800
$ws.Range("K136").Value = 6031.0002
$ws.Range("L136").Value = 11400
$ws.Range("M136").Value = -3481.0002
$ws.Range("N136").Value = -16500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 13161005
$ws.Range("I122").Value = 35716450
$ws.Range("J122").Value = 3662.9167
$ws.Range("K122").Value = 107149350
$ws.Range("L122").Value = 10988.7501
$ws.Range("M122").Value = -107146900
$ws.Range("N122").Value = -15888.7501
$ws.Range("H123").Value = 22426.387
$ws.Range("J123").Value = 22426.387
$ws.Range("L123").Value = 22426.387
$ws.Range("N123").Value = -32226.387
$ws.Range("H125").Value = 60182
$ws.Range("J125").Value = 60182
$ws.Range("L125").Value = 60182
$ws.Range("N125").Value = -70022
$ws.Range("H131").Value = 60536
$ws.Range("J131").Value = 60536
$ws.Range("L131").Value = 60536
$ws.Range("N131").Value = -70616
$ws.Range("H132").Value = 2915.9033
$ws.Range("I132").Value = 2240.5881
$ws.Range("J132").Value = 3735.9285
$ws.Range("K132").Value = 6721.7643
$ws.Range("L132").Value = 11207.7855
$ws.Range("M132").Value = -4191.7643
$ws.Range("N132").Value = -16267.7855
$ws.Range("H136").Value = 2480.2415
$ws.Range("I136").Value = 2048.5789
$ws.Range("J136").Value = 3300.4
$ws.Range("K136").Value = 6145.736699999999
$ws.Range("L136").Value = 9901.200000000001
$ws.Range("M136").Value = -3595.736699999999
$ws.Range("N136").Value = -15001.2
